# Apply "Trade #63 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - refresh aggregate metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.12   # Current Capital
$summary.Range("B4").Value = 0.11      # Total P&L $
$summary.Range("B6").Value = 63        # Total Trades
$summary.Range("B7").Value = 28        # Winning Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4) refresh
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.12     # Capital
$status.Range("D4").Value = 63         # Trades
$status.Range("E4").Value = 0.11       # P&L $
$status.Range("F4").Value = 0.12       # P&L %
$status.Range("G4").Value = 44.44      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new closed trade (row 64) to both "All Trades" and
#    "MarketMaking" sheets - they carry the same trade log.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A64").Value = 63
    # Force text format on the date cell so Excel does not auto-convert
    # the "2026-02-17" literal into a date serial value.
    $ws.Range("B64").NumberFormat = "@"
    $ws.Range("B64").Value = "2026-02-17"
    $ws.Range("C64").Value = "12:53:13"
    $ws.Range("D64").Value = "MarketMaking"
    $ws.Range("E64").Value = "DOWN"
    $ws.Range("F64").Value = 0.86
    $ws.Range("G64").Value = 0.87
    $ws.Range("H64").Value = "CLOSED"
    $ws.Range("I64").Value = 1.1628
    $ws.Range("J64").Value = 0.01
    $ws.Range("K64").Value = 100.12
    $ws.Range("L64").Value = 0
    $ws.Range("M64").Value = 0
    $ws.Range("N64").Value = 0.6
    $ws.Range("O64").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P64").Value = "early_exit"
    $ws.Range("Q64").Value = 0.14
}
